$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.554.06"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "1.666.15"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "238.38"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.4802"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.2637"
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").Value = "0.06178"
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("D10").Value = "0.07112"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").Value = "1.660.91"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "14.84"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "0.5900"
$ws.Range("E13").Value = "  -5.07%  "
$ws.Range("D14").Value = "4.376"
$ws.Range("E14").Value = "  -5.07%  "
$ws.Range("D15").Value = "75.23"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "0.9996"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "25.543.06"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "0.000006778"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("D20").Value = "11.48"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "1.874.15"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "4.425"
$ws.Range("E22").Value = "  -3.45%  "
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "5.289"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "135.30"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("D26").Value = "15.04"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "1.382"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "105.05"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").Value = "1.715"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").Value = "3.966"
$ws.Range("E30").Value = "  +4.95%  "
$ws.Range("D31").Value = "3.655"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "0.07734"
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "0.04230"
$ws.Range("E34").Value = "  -8.51%  "
$ws.Range("D35").Value = "2.599"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "0.6142"
$ws.Range("E36").Value = "  +6.31%  "
$ws.Range("D37").Value = "0.9539"
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("D38").Value = "2.590"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").Value = "0.8612"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").Value = "0.9994"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "1.854"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "0.01469"
$ws.Range("E42").Value = "  -6.15%  "
$ws.Range("D43").Value = "97.35"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "0.3765"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "4.850"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").Value = "0.1123"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").Value = "6.233"
$ws.Range("D48").Value = "0.05262"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D51").Value = "7.341"
$ws.Range("E51").Value = "  +1.77%  "
